# Updates the cryptos list (prices / 1h volume %, and a few re-ranked rows)
# to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price cells hold values that look numeric (e.g. '245.00', '8.70').
# They are stored as plain text in the sheet (to keep trailing zeros and
# the site's dotted-thousands formatting like '43.101.49'), so force those
# specific cells to Text format first; otherwise Excel would silently
# coerce them into numbers and drop the meaningful trailing zero.
$textCells = @('D5', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D15', 'D20', 'D22', 'D23', 'D24', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range('D2').Value = '43.101.49'
$ws.Range('E2').Value = '  +4.54%  '
$ws.Range('D3').Value = '2.247.09'
$ws.Range('E3').Value = '  +3.82%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '245.00'
$ws.Range('E5').Value = '  +3.73%  '
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('D7').Value = '75.29'
$ws.Range('E7').Value = '  +8.52%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '0.606'
$ws.Range('D10').Value = '41.22'
$ws.Range('E10').Value = '  +6.56%  '
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  +2.65%  '
$ws.Range('D12').Value = '6.94'
$ws.Range('D13').Value = '0.101'
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = '2.585.40'
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').Value = '14.59'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').Value = '2.234.72'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('D18').Value = '43.005.88'
$ws.Range('E18').Value = '  +4.89%  '
$ws.Range('E19').Value = '  +5.91%  '
$ws.Range('D20').Value = '71.05'
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('E21').Value = '  +3.53%  '
$ws.Range('D22').Value = '9.90'
$ws.Range('E22').Value = '  +6.39%  '
$ws.Range('D23').Value = '229.72'
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').Value = '2.19'
$ws.Range('E24').Value = '  +16.94%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '10.87'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').Value = '3.40'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('E28').Value = '  +2.68%  '
$ws.Range('D29').Value = '38.66'
$ws.Range('E29').Value = '  +29.21%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.22'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '172.19'
$ws.Range('E31').Value = '  +2.14%  '
$ws.Range('D32').Value = '20.29'
$ws.Range('E32').Value = '  +2.63%  '
$ws.Range('D33').Value = '0.0796'
$ws.Range('E33').Value = '  +5.84%  '
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('E35').Value = '  +2.22%  '
$ws.Range('E36').Value = '  +7.75%  '
$ws.Range('D37').Value = '4.36'
$ws.Range('E37').Value = '  +7.09%  '
$ws.Range('D38').Value = '0.0332'
$ws.Range('E38').Value = '  +19.34%  '
$ws.Range('D39').Value = '13.14'
$ws.Range('E39').Value = '  +14.95%  '
$ws.Range('E40').Value = '  +3.66%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.206'
$ws.Range('E41').Value = '  +10.63%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').Value = '5.46'
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('D43').Value = '59.44'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('D44').Value = '104.92'
$ws.Range('E44').Value = '  +8.72%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '8.70'
$ws.Range('E45').Value = '  +5.98%  '
$ws.Range('B46').Value = 'WOONetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D46').Value = '0.488'
$ws.Range('E46').Value = '  +32.56%  '
$ws.Range('E47').Value = '  +3.39%  '
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  +11.77%  '
$ws.Range('E49').Value = '  +3.64%  '
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('D51').Value = '2.457.98'
$ws.Range('E51').Value = '  +4.20%  '
